# Insert a new data row before current row 8 (HEALSEC row), pushing the
# existing rows (and the totals/footer rows) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(8).Insert()

# Copy the cell formatting (styles) from the row that used to be row 8
# (now shifted to row 9) onto the newly inserted blank row 8, so the new
# row visually matches the rest of the table instead of getting Excel's
# generic "inserted row" styling.
$ws.Range("A9:N9").Copy()
$ws.Range("A8:N8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows(8).RowHeight = 25.5

# Re-create the merged cells for the new row, matching the pattern used by
# every other product row in the table.
$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()

# Fill in the data for the new medicine row.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "DOSTINEX 0.5 MG 2 TABS."
$ws.Range("H8").Value = "0:0"
$ws.Range("L8").Value = 172
$ws.Range("N8").Value = "1:0"

# The "م" (sequence number) column is a plain typed sequence, not a
# formula, and keeps counting 1..17 downward regardless of which product
# ended up in which row, so re-stamp it for every row now that a row was
# inserted.
for ($r = 8; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}

# Update the running total to include the newly added row.
$ws.Range("K21").Value = 722
